$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 was a plain text label "2023-06-03"; turn it into a real date value
# formatted as yyyy-mm-dd hh:mm:ss (same serial date, new display format).
$ws.Range("B2:B4").NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Range("B2").Value = "2023-06-03"

# G3/G4 previously carried a special (9pt) note style; drop back to the
# sheet's normal/default style.
$ws.Range("G3:G4").Style = "常规"

# Update the note text in G4.
$ws.Range("G4").Value = "天神下凡酣畅淋漓，她腿软了，我也舒服了，好热好热"

# Selection moved to A2:A4 with A2 active.
$ws.Range("A2:A4").Select()
